$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the FECHA (date/time) cell with the new timestamp
$ws.Range("C6").Value = "16/11/2024 17:03"

# 2) Update attendance marks in column S (day 16)
#    Style indices reused from the workbook:
#      s=34 -> "P" (green)   donor: S8  (stays P)
#      s=35 -> "F" (red)     donor: S11 (stays F)
#      s=36 -> "R" (yellow)  donor: S17 (currently R, changes later - copy first!)
#      s=26 -> empty numeric donor: S33 (stays empty)

# Capture the "R" (yellow) formatting BEFORE S17 itself is changed, and apply
# it to the rows that need to become "R".
$ws.Range("S17").Copy()
$ws.Range("S15").PasteSpecial(-4122)
$ws.Range("S17").Copy()
$ws.Range("S18").PasteSpecial(-4122)
$ws.Range("S15").Value = "R"
$ws.Range("S18").Value = "R"

# Now change S17 itself from "R" to "F"
$ws.Range("S11").Copy()
$ws.Range("S17").PasteSpecial(-4122)
$ws.Range("S17").Value = "F"

# Rows 9, 10, 14 go from "P" to "F"
$ws.Range("S11").Copy()
$ws.Range("S9").PasteSpecial(-4122)
$ws.Range("S11").Copy()
$ws.Range("S10").PasteSpecial(-4122)
$ws.Range("S11").Copy()
$ws.Range("S14").PasteSpecial(-4122)
$ws.Range("S9").Value = "F"
$ws.Range("S10").Value = "F"
$ws.Range("S14").Value = "F"

# Rows 19-32 get cleared back to plain empty numeric cells
$emptyRows = 19..32
foreach ($r in $emptyRows) {
    $ws.Range("S33").Copy()
    $ws.Range("S$r").PasteSpecial(-4122)
    $ws.Range("S$r").ClearContents()
}

$excel.CutCopyMode = 0
